{"js": "// Completa la actividad 1.1 y 1.2: rellena los datos del estudiante\n// (Nombre, Rut, Carrera, Sede) en la tabla de \"Antecedentes Personales\".\n//\n// La tabla es la segunda tabla del documento (\u00edndice 1, base 0):\n//   Fila 0: Nombre estudiante | <vac\u00edo>\n//   Fila 1: Rut               | <vac\u00edo>\n//   Fila 2: Carrera           | <vac\u00edo>\n//   Fila 3: Sede              | <vac\u00edo>\n// Cada celda vac\u00eda de la columna 2 ya tiene un p\u00e1rrafo con negrita\n// definida a nivel de marca de p\u00e1rrafo; solo hace falta insertar el\n// texto (en negrita) al final de esa celda.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst studentTable = tables.items[1];\n\nasync function setCellValue(rowIndex, text) {\n  const cell = studentTable.getCell(rowIndex, 1);\n  const range = cell.body.insertText(text, \"End\");\n  range.font.bold = true;\n  await context.sync();\n}\n\nawait setCellValue(0, \"David Andr\u00e9s Zurita Corval\u00e1n\");\nawait setCellValue(1, \"17475735-6\");\nawait setCellValue(2, \"Ingenier\u00eda en Inform\u00e1tica\");\nawait setCellValue(3, \"Vina del Mar\");\n", "ps1": "# Completa la actividad 1.1 y 1.2: rellena los datos del estudiante\n# (Nombre, Rut, Carrera, Sede) en la tabla de \"Antecedentes Personales\".\n#\n# La tabla es la segunda tabla del documento:\n#   Fila 1: Nombre estudiante | <vac\u00edo>\n#   Fila 2: Rut               | <vac\u00edo>\n#   Fila 3: Carrera           | <vac\u00edo>\n#   Fila 4: Sede              | <vac\u00edo>\n# Cada celda vac\u00eda de la columna 2 ya tiene un p\u00e1rrafo con rPr/b definido\n# (el valor se muestra en negrita); solo hace falta insertar el texto.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(2)\n\nfunction Set-CellValue($table, $row, $col, $text) {\n    $cell = $table.Cell($row, $col)\n    $paraRange = $cell.Range.Paragraphs.Item(1).Range\n    $startPos = $paraRange.Start\n    $paraRange.InsertBefore($text)\n    $valueRange = $d.Range($startPos, $startPos + $text.Length)\n    $valueRange.Font.Bold = $true\n}\n\nSet-CellValue $table 1 2 \"David Andr\u00e9s Zurita Corval\u00e1n\"\nSet-CellValue $table 2 2 \"17475735-6\"\nSet-CellValue $table 3 2 \"Ingenier\u00eda en Inform\u00e1tica\"\nSet-CellValue $table 4 2 \"Vina del Mar\"\n"}
